$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4120.6
$ws.Range("J29").Value = 9899
$ws.Range("L29").Value = 29697
$ws.Range("N29").Value = -30259
$ws.Range("H31").Value = 4525
$ws.Range("I31").Value = 4525
$ws.Range("K31").Value = 13575
$ws.Range("M31").Value = -13345
$ws.Range("H42").Value = 1581.6666
$ws.Range("I42").Value = 501.66666
$ws.Range("K42").Value = 1504.99998
$ws.Range("M42").Value = -1274.99998
$ws.Range("H46").Value = 9709
$ws.Range("J46").Value = 7418
$ws.Range("L46").Value = 22254
$ws.Range("N46").Value = -22492
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = ""
$ws.Range("N49").Value = 0
$ws.Range("H60").Value = 9709
$ws.Range("J60").Value = 7418
$ws.Range("L60").Value = 22254
$ws.Range("N60").Value = -23222
$ws.Range("H62").Value = 2499.5
$ws.Range("J62").Value = 2499
$ws.Range("L62").Value = 2499
$ws.Range("N62").Value = -3747
$ws.Range("H65").Value = 2499.5
$ws.Range("J65").Value = 2499
$ws.Range("L65").Value = 12495
$ws.Range("N65").Value = -18735
$ws.Range("H137").Value = 1356
$ws.Range("I137").Value = 1356
$ws.Range("K137").Value = 4068
$ws.Range("M137").Value = -1518

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 579.7778
$ws.Range("I2").Value = 174.14285
$ws.Range("K2").Value = 174.14285
$ws.Range("M2").Value = -61.14285000000001
$ws.Range("I45").Value = 3999
$ws.Range("K45").Value = 3999
$ws.Range("M45").Value = -3622
$ws.Range("H61").Value = 2666
$ws.Range("I61").Value = 2666
$ws.Range("K61").Value = 2666
$ws.Range("M61").Value = -2454
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = ""
$ws.Range("N109").Value = 0
$ws.Range("H116").Value = 579.7778
$ws.Range("I116").Value = 174.14285
$ws.Range("K116").Value = 174.14285
$ws.Range("M116").Value = 2119.85715
$ws.Range("H136").Value = 2666
$ws.Range("I136").Value = 2666
$ws.Range("K136").Value = 7998
$ws.Range("M136").Value = -5448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 579.7778
$ws.Range("I3").Value = 174.14285
$ws.Range("K3").Value = 174.14285
$ws.Range("M3").Value = -60.14285000000001
$ws.Range("H134").Value = 4011
$ws.Range("I134").Value = 4011
$ws.Range("K134").Value = 12033
$ws.Range("M134").Value = -9498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = ""
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = ""
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = 0
$ws.Range("H35").Value = 277.5
$ws.Range("I35").Value = 277.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 277.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = 16.5
$ws.Range("H58").Value = 4772.8
$ws.Range("I58").Value = 1941.8572
$ws.Range("J58").Value = 7249.875
$ws.Range("K58").Value = 1941.8572
$ws.Range("L58").Value = 7249.875
$ws.Range("M58").Value = -1738.8572
$ws.Range("N58").Value = -7655.875
$ws.Range("H136").Value = 4772.8
$ws.Range("I136").Value = 1941.8572
$ws.Range("J136").Value = 7249.875
$ws.Range("K136").Value = 5825.571599999999
$ws.Range("L136").Value = 21749.625
$ws.Range("M136").Value = -3275.571599999999
$ws.Range("N136").Value = -26849.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 3000
$ws.Range("J123").Value = 3000
$ws.Range("L123").Value = 9000
$ws.Range("N123").Value = -13900
$ws.Range("H140").Value = 5000
$ws.Range("I140").Value = 5000
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 15000
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = ""
$ws.Range("N140").Value = -9820
$ws.Range("H141").Value = 2989
$ws.Range("I141").Value = 2989
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8967
$ws.Range("L141").Value = ""
$ws.Range("N141").Value = 0
$ws.Range("M141").Value = -3787

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = ""
$ws.Range("N126").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3709.1428
$ws.Range("I122").Value = 3192.8
$ws.Range("K122").Value = 9578.400000000001
$ws.Range("M122").Value = -7128.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = ""
$ws.Range("N82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = ""
$ws.Range("N85").Value = 0
$ws.Range("H107").Value = 1393
$ws.Range("J107").Value = 1800.5
$ws.Range("L107").Value = 5401.5
$ws.Range("N107").Value = -9241.5
$ws.Range("H117").Value = 33499.5
$ws.Range("J117").Value = 33499.5
$ws.Range("L117").Value = 33499.5
$ws.Range("N117").Value = -42677.5
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = ""
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = 0
$ws.Range("H132").Value = 1490
$ws.Range("I132").Value = 1548
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 4644
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -2114
$ws.Range("N132").Value = -8660
$ws.Range("H136").Value = 3454.4443
$ws.Range("I136").Value = 1584.7142
$ws.Range("J136").Value = 9998.5
$ws.Range("K136").Value = 4754.142599999999
$ws.Range("L136").Value = 29995.5
$ws.Range("M136").Value = -2204.142599999999
$ws.Range("N136").Value = -35095.5
